# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> clrScheme name="Office"     (name="Office Theme")
#   ppt/theme/theme2.xml  -> clrScheme name="Red Violet"  (name="Integral")
# theme2.xml is the theme actually wired to the (single) slide master /
# the slides, so it is the one the Color Scheme UI (and this COM surface)
# can reach. The authored edit swaps the two themes' contents, which -
# for the part that actually drives what renders on the slides - means
# the "Red Violet" palette used by theme2.xml is replaced by the
# "Office" palette that used to live in theme1.xml.
#
# Slide.ThemeColorScheme exposes all twelve theme colour slots (dk1, lt1,
# dk2, lt2, accent1-6, hlink, folHlink) in clrScheme document order and
# writing .RGB on them edits the underlying theme XML in place (unlike
# the legacy 8-slot Slide.ColorScheme, which blanks the clrScheme name).
# Because there is a single master/theme for the whole deck, this change
# is global - any slide can be used as the anchor.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      -> 000000
$tcs.Item(2).RGB  = 16777215   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      -> 44546A
$tcs.Item(4).RGB  = 15132391   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  -> ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  -> FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  -> 4472C4
$tcs.Item(10).RGB = 4697456    # accent6  -> 70AD47
$tcs.Item(11).RGB = 12673797   # hlink    -> 0563C1
$tcs.Item(12).RGB = 7491477    # folHlink -> 954F72
